$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert two new rows right after the header (before THIAGO row) ---
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2,1).Value = "'004572740"
$ws.Cells.Item(2,2).Value = "PAULO"
$ws.Cells.Item(2,3).Value = 231567.92

$ws.Cells.Item(3,1).Value = "'005726697"
$ws.Cells.Item(3,2).Value = "FERNANDO"
$ws.Cells.Item(3,3).Value = 120000

# Row 4 is now THIAGO (unchanged)

# --- Step 2: insert nine new rows right before CAIO (currently row 5) ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5,1).Value = "'004001621"
$ws.Cells.Item(5,2).Value = "DANIELA"
$ws.Cells.Item(5,3).Value = 20000

$ws.Cells.Item(6,1).Value = "'004452476"
$ws.Cells.Item(6,2).Value = "IVONE"
$ws.Cells.Item(6,3).Value = 19937.62

$ws.Cells.Item(7,1).Value = "'005305448"
$ws.Cells.Item(7,2).Value = "ALPHASITIO"
$ws.Cells.Item(7,3).Value = 916.92

$ws.Cells.Item(8,1).Value = "'004724018"
$ws.Cells.Item(8,2).Value = "ASPA"
$ws.Cells.Item(8,3).Value = 910.58

$ws.Cells.Item(9,1).Value = "'004392159"
$ws.Cells.Item(9,2).Value = "RODRIGO"
$ws.Cells.Item(9,3).Value = 900.21

$ws.Cells.Item(10,1).Value = "'004870019"
$ws.Cells.Item(10,2).Value = "MARIA"
$ws.Cells.Item(10,3).Value = 760.19

$ws.Cells.Item(11,1).Value = "'005685353"
$ws.Cells.Item(11,2).Value = "CARLOS"
$ws.Cells.Item(11,3).Value = 550.73

$ws.Cells.Item(12,1).Value = "'004862672"
$ws.Cells.Item(12,2).Value = "RENATO"
$ws.Cells.Item(12,3).Value = 526.58

$ws.Cells.Item(13,1).Value = "'004713953"
$ws.Cells.Item(13,2).Value = "ALESSANDRA"
$ws.Cells.Item(13,3).Value = 481.75

# Row 14 is now CAIO - update its Saldo value
$ws.Cells.Item(14,3).Value = 473.12

# --- Step 3: remove the old duplicate rows that used to follow CAIO ---
# They now sit at rows 15-21 (ALPHASITIO, ASPA, RODRIGO, MARIA, CARLOS, RENATO, ALESSANDRA)
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(15).Delete()
